$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "41.777.82"
Set-TextValue $ws.Range("E2") "  +4.86%  "

Set-TextValue $ws.Range("D3") "2.273.21"
Set-TextValue $ws.Range("E3") "  +3.71%  "

Set-TextValue $ws.Range("E4") "  -0.01%  "

Set-TextValue $ws.Range("D5") "302.88"
Set-TextValue $ws.Range("E5") "  +3.67%  "

Set-TextValue $ws.Range("D6") "93.10"
Set-TextValue $ws.Range("E6") "  +7.38%  "

Set-TextValue $ws.Range("E7") "  +2.98%  "

Set-TextValue $ws.Range("E8") "  -0.06%  "

Set-TextValue $ws.Range("D9") "0.489"
Set-TextValue $ws.Range("E9") "  +5.28%  "

Set-TextValue $ws.Range("D10") "54.47"
Set-TextValue $ws.Range("E10") "  +8.86%  "

Set-TextValue $ws.Range("D11") "32.36"
Set-TextValue $ws.Range("E11") "  +8.13%  "

Set-TextValue $ws.Range("D12") "0.0799"
Set-TextValue $ws.Range("E12") "  +3.08%  "

Set-TextValue $ws.Range("E13") "  +3.47%  "

Set-TextValue $ws.Range("D14") "6.67"
Set-TextValue $ws.Range("E14") "  +4.06%  "

Set-TextValue $ws.Range("D15") "2.618.11"
Set-TextValue $ws.Range("E15") "  +3.56%  "

Set-TextValue $ws.Range("D16") "14.23"
Set-TextValue $ws.Range("E16") "  +4.44%  "

Set-TextValue $ws.Range("D17") "2.234.49"
Set-TextValue $ws.Range("E17") "  -0.79%  "

Set-TextValue $ws.Range("D18") "0.757"
Set-TextValue $ws.Range("E18") "  +4.72%  "

Set-TextValue $ws.Range("D19") "41.626.82"
Set-TextValue $ws.Range("E19") "  +4.85%  "

Set-TextValue $ws.Range("D20") "12.50"
Set-TextValue $ws.Range("E20") "  +12.37%  "

Set-TextValue $ws.Range("D21") "0.0₃0911"
Set-TextValue $ws.Range("E21") "  +3.52%  "

Set-TextValue $ws.Range("D23") "67.15"
Set-TextValue $ws.Range("E23") "  +3.28%  "

Set-TextValue $ws.Range("D24") "240.76"
Set-TextValue $ws.Range("E24") "  +2.00%  "

Set-TextValue $ws.Range("E25") "  +5.99%  "

Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.13%  "

Set-TextValue $ws.Range("D27") "1.88"
Set-TextValue $ws.Range("E27") "  +4.60%  "

Set-TextValue $ws.Range("D28") "23.86"
Set-TextValue $ws.Range("E28") "  +4.18%  "

Set-TextValue $ws.Range("E29") "  +6.62%  "

Set-TextValue $ws.Range("D30") "9.73"
Set-TextValue $ws.Range("E30") "  +6.73%  "

Set-TextValue $ws.Range("D31") "34.31"
Set-TextValue $ws.Range("E31") "  +9.76%  "

Set-TextValue $ws.Range("D32") "157.69"
Set-TextValue $ws.Range("E32") "  +1.26%  "

Set-TextValue $ws.Range("D33") "1.00"
Set-TextValue $ws.Range("E33") "  +0.11%  "

Set-TextValue $ws.Range("E34") "  +6.82%  "

Set-TextValue $ws.Range("B35") "LidoDAOToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D35") "3.08"
Set-TextValue $ws.Range("E35") "  +9.95%  "

Set-TextValue $ws.Range("B36") "Hedera"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.0739"
Set-TextValue $ws.Range("E36") "  +4.74%  "

Set-TextValue $ws.Range("E37") "  +1.87%  "

Set-TextValue $ws.Range("D38") "16.76"
Set-TextValue $ws.Range("E38") "  +11.12%  "

Set-TextValue $ws.Range("E39") "  +7.58%  "

Set-TextValue $ws.Range("E40") "  +3.07%  "

Set-TextValue $ws.Range("D41") "1.79"
Set-TextValue $ws.Range("E41") "  +7.46%  "

Set-TextValue $ws.Range("D42") "4.01"
Set-TextValue $ws.Range("E42") "  +7.55%  "

Set-TextValue $ws.Range("D43") "20.39"
Set-TextValue $ws.Range("E43") "  +18.92%  "

Set-TextValue $ws.Range("D44") "2.062.42"
Set-TextValue $ws.Range("E44") "  -2.79%  "

Set-TextValue $ws.Range("D45") "0.0279"
Set-TextValue $ws.Range("E45") "  +4.71%  "

Set-TextValue $ws.Range("B46") "NEARProtocol"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D46") "2.98"
Set-TextValue $ws.Range("E46") "  +12.64%  "

Set-TextValue $ws.Range("B47") "FraxShare"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D47") "10.10"
Set-TextValue $ws.Range("E47") "  +4.50%  "

Set-TextValue $ws.Range("D48") "1.97"
Set-TextValue $ws.Range("E48") "  -5.81%  "

Set-TextValue $ws.Range("D49") "2.489.71"
Set-TextValue $ws.Range("E49") "  +3.79%  "

Set-TextValue $ws.Range("E50") "  +3.82%  "

Set-TextValue $ws.Range("E51") "  +5.10%  "
